$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1800766283524904
$ws.Range("C2").Value = 0.5938697318007663
$ws.Range("J2").Value = 0.007662835249042145
$ws.Range("P2").Value = 0.1647509578544061
$ws.Range("S2").Value = 0.05363984674329502
$ws.Range("B3").Value = 0.01273885350318471
$ws.Range("C3").Value = 0.01910828025477707
$ws.Range("J3").Value = 0.006369426751592357
$ws.Range("P3").Value = 0.7579617834394905
$ws.Range("S3").Value = 0.2038216560509554
$ws.Range("J4").Value = 0.06521739130434782
$ws.Range("P4").Value = 0.6956521739130435
$ws.Range("S4").Value = 0.2391304347826087
$ws.Range("B6").Value = 0.05911330049261083
$ws.Range("D6").Value = 0.01970443349753695
$ws.Range("F6").Value = 0.07389162561576355
$ws.Range("J6").Value = 0.2463054187192118
$ws.Range("O6").Value = 0.02463054187192118
$ws.Range("Q6").Value = 0.1477832512315271
$ws.Range("R6").Value = 0.0541871921182266
$ws.Range("S6").Value = 0.374384236453202
$ws.Range("B7").Value = 0.07870370370370371
$ws.Range("D7").Value = 0.01388888888888889
$ws.Range("F7").Value = 0.02777777777777778
$ws.Range("J7").Value = 0.1527777777777778
$ws.Range("O7").Value = 0.009259259259259259
$ws.Range("Q7").Value = 0.2083333333333333
$ws.Range("R7").Value = 0.07870370370370371
$ws.Range("S7").Value = 0.4305555555555556
$ws.Range("B8").Value = 0.07610993657505286
$ws.Range("D8").Value = 0.01902748414376321
$ws.Range("F8").Value = 0.07188160676532769
$ws.Range("J8").Value = 0.1353065539112051
$ws.Range("O8").Value = 0.008456659619450317
$ws.Range("Q8").Value = 0.1818181818181818
$ws.Range("R8").Value = 0.105708245243129
$ws.Range("S8").Value = 0.4016913319238901
$ws.Range("B9").Value = 0.09195402298850575
$ws.Range("D9").Value = 0.02873563218390805
$ws.Range("F9").Value = 0.05172413793103448
$ws.Range("J9").Value = 0.132183908045977
$ws.Range("O9").Value = 0.01149425287356322
$ws.Range("Q9").Value = 0.2413793103448276
$ws.Range("R9").Value = 0.05747126436781609
$ws.Range("S9").Value = 0.3850574712643678
$ws.Range("B10").Value = 0.1090764331210191
$ws.Range("D10").Value = 0.01990445859872611
$ws.Range("F10").Value = 0.05812101910828026
$ws.Range("J10").Value = 0.1170382165605096
$ws.Range("O10").Value = 0.01592356687898089
$ws.Range("Q10").Value = 0.2197452229299363
$ws.Range("R10").Value = 0.08200636942675159
$ws.Range("S10").Value = 0.3781847133757962
$ws.Range("G11").Value = 0.1656804733727811
$ws.Range("J11").Value = 0.07692307692307693
$ws.Range("K11").Value = 0.2218934911242604
$ws.Range("L11").Value = 0.514792899408284
$ws.Range("S11").Value = 0.02071005917159763
$ws.Range("G12").Value = 0.7329842931937173
$ws.Range("J12").Value = 0.2094240837696335
$ws.Range("L12").Value = 0.03141361256544502
$ws.Range("S12").Value = 0.02617801047120419
$ws.Range("G13").Value = 0.6938775510204082
$ws.Range("J13").Value = 0.2653061224489796
$ws.Range("S13").Value = 0.04081632653061224
$ws.Range("F15").Value = 0.01041666666666667
$ws.Range("H15").Value = 0.1927083333333333
$ws.Range("I15").Value = 0.07291666666666667
$ws.Range("J15").Value = 0.375
$ws.Range("K15").Value = 0.05208333333333334
$ws.Range("M15").Value = 0.02083333333333333
$ws.Range("O15").Value = 0.05208333333333334
$ws.Range("S15").Value = 0.2239583333333333
$ws.Range("F16").Value = 0.01612903225806452
$ws.Range("H16").Value = 0.1774193548387097
$ws.Range("I16").Value = 0.03763440860215054
$ws.Range("J16").Value = 0.4032258064516129
$ws.Range("K16").Value = 0.1397849462365591
$ws.Range("M16").Value = 0.04301075268817205
$ws.Range("O16").Value = 0.06451612903225806
$ws.Range("S16").Value = 0.1182795698924731
$ws.Range("F17").Value = 0.02306079664570231
$ws.Range("H17").Value = 0.1593291404612159
$ws.Range("I17").Value = 0.09433962264150944
$ws.Range("J17").Value = 0.4339622641509434
$ws.Range("K17").Value = 0.09853249475890985
$ws.Range("M17").Value = 0.01677148846960168
$ws.Range("O17").Value = 0.06289308176100629
$ws.Range("S17").Value = 0.1111111111111111
$ws.Range("F18").Value = 0.02094240837696335
$ws.Range("H18").Value = 0.2198952879581152
$ws.Range("I18").Value = 0.06806282722513089
$ws.Range("J18").Value = 0.4031413612565445
$ws.Range("K18").Value = 0.08900523560209424
$ws.Range("M18").Value = 0.02094240837696335
$ws.Range("O18").Value = 0.05759162303664921
$ws.Range("S18").Value = 0.1204188481675393
$ws.Range("F19").Value = 0.01892744479495268
$ws.Range("H19").Value = 0.2271293375394322
$ws.Range("I19").Value = 0.07097791798107256
$ws.Range("J19").Value = 0.3501577287066246
$ws.Range("K19").Value = 0.1190851735015773
$ws.Range("M19").Value = 0.01735015772870662
$ws.Range("O19").Value = 0.05520504731861198
$ws.Range("S19").Value = 0.1411671924290221
